$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write a numeric (date) value into A3
$ws.Range("A3").Value = 43101
$ws.Range("A3").NumberFormat = "d-mmm-yy"

# Auto-fit column A width to match the newly written content ("bestFit").
# (9.5 "characters" is the input that rounds to the engine's nearest
# representable column width to Excel's real bestFit width for "1-Jan-18".)
$ws.Columns.Item(1).ColumnWidth = 9.5

# Update the active selection to E8, matching the target workbook state
$ws.Range("E8").Select()
